$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.828.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.053.27"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.53"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.38%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +5.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.65"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.55%  "

$ws.Range("E10").Value = "  +7.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.373"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.65%  "

$ws.Range("E12").Value = "  +2.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.578.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.76%  "

$ws.Range("E15").Value = "  +16.65%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.773.63"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.38%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.044.18"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.92%  "

$ws.Range("E20").Value = "  +5.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "338.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.88%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.85%  "

$ws.Range("E25").Value = "  +6.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0979"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.74%  "

$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.39"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +10.55%  "

$ws.Range("E30").Value = "  +6.73%  "

$ws.Range("E31").Value = "  +5.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.98%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.75"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.16%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.55"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.70%  "

$ws.Range("E35").Value = "  +7.16%  "

$ws.Range("E36").Value = "  +3.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.17"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0703"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.090.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.41%  "

$ws.Range("E40").Value = "  +3.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.88"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +8.65%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.26%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.664"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.329.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.22%  "

$ws.Range("E46").Value = "  +3.59%  "

$ws.Range("E47").Value = "  +3.17%  "

$ws.Range("E48").Value = "  +3.78%  "

$ws.Range("E49").Value = "  +4.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0901"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.44%  "
